$wb = $excel.ActiveWorkbook

# Add a new worksheet "5.1" after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "5.1"

# Column A labels first (matches shared-string insertion order: Country, Mexico, China, India, Phillippines, Dominican Republic)
$ws.Range("A1").Value = "Country"
$ws.Range("A2").Value = "Mexico"
$ws.Range("A3").Value = "China"
$ws.Range("A4").Value = "India"
$ws.Range("A5").Value = "Phillippines"
$ws.Range("A6").Value = "Dominican Republic"

# Header row (CA, NY, FL, TX, NJ, Totals)
$ws.Range("B1").Value = "CA"
$ws.Range("C1").Value = "NY"
$ws.Range("D1").Value = "FL"
$ws.Range("E1").Value = "TX"
$ws.Range("F1").Value = "NJ"
$ws.Range("G1").Value = "Totals"

# Totals row label (reuses the "Totals" shared string created above)
$ws.Range("A7").Value = "Totals"

# Data rows
$ws.Range("B2").Value = 50645
$ws.Range("C2").Value = 2437
$ws.Range("D2").Value = 3113
$ws.Range("E2").Value = 32811
$ws.Range("F2").Value = 2437
$ws.Range("G2").Value = 139120

$ws.Range("B3").Value = 18680
$ws.Range("C3").Value = 18859
$ws.Range("D3").Value = 1620
$ws.Range("E3").Value = 3280
$ws.Range("F3").Value = 2253
$ws.Range("G3").Value = 70863

$ws.Range("B4").Value = 15099
$ws.Range("C4").Value = 5116
$ws.Range("D4").Value = 2019
$ws.Range("E4").Value = 5777
$ws.Range("F4").Value = 8123
$ws.Range("G4").Value = 69162

$ws.Range("B5").Value = 24082
$ws.Range("C5").Value = 2361
$ws.Range("D5").Value = 2320
$ws.Range("E5").Value = 2525
$ws.Range("F5").Value = 2321
$ws.Range("G5").Value = 58173

$ws.Range("B6").Value = 172
$ws.Range("C6").Value = 26249
$ws.Range("D6").Value = 3900
$ws.Range("E6").Value = 275
$ws.Range("F6").Value = 8444
$ws.Range("G6").Value = 53870

$ws.Range("B7").Value = 208446
$ws.Range("C7").Value = 147999
$ws.Range("D7").Value = 107276
$ws.Range("E7").Value = 87750
$ws.Range("F7").Value = 56920
$ws.Range("G7").Value = 1042625
